$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23: update title (D23) and link (E23)
$ws.Range("D23").Value = "딥러닝 기본과 NLP를 익히는데 도움이 될 만한 최신 (2020년 2021년) 동영상 강좌 13종입니다.`n하나 하나 직접 들어본 분의 추천이니"
$ws.Range("E23").Value = "https://theonly1.tistory.com/2764"

# Row 39: update title (D39) and link (E39)
$ws.Range("D39").Value = "Getting to know probability distributions"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/Getting-to-know-probability-distributions-1"

# Row 51: update title (D51) and link (E51)
$ws.Range("D51").Value = "삼성 가성비 복합기 SL-J1660 사용 심플 후기 (드라이버 다운로드)"
$ws.Range("E51").Value = "https://bskyvision.com/1159"
